$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 98, shifting existing rows 98:106 down to 99:107.
$ws.Rows.Item(98).Insert()

# Populate the new row 98 with the new weekly record.
$ws.Cells.Item(98, 1).Value = 10
$ws.Cells.Item(98, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(98, 3).Value = "La Araucanía"
$ws.Cells.Item(98, 4).Value = 44776
$ws.Cells.Item(98, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(98, 5).Value = 9
$ws.Cells.Item(98, 6).Value = 100112035
$ws.Cells.Item(98, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(98, 8).Value = "Sin especificar"
$ws.Cells.Item(98, 9).Value = "Primera"
$ws.Cells.Item(98, 10).Value = 50
$ws.Cells.Item(98, 11).Value = 25000
$ws.Cells.Item(98, 12).Value = 25000
$ws.Cells.Item(98, 13).Value = 25000
$ws.Cells.Item(98, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(98, 15).Value = "Región Metropolitana"
$ws.Cells.Item(98, 16).Value = 2500
$ws.Cells.Item(98, 17).Value = 10
$ws.Cells.Item(98, 18).Value = "Hortaliza"
